$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# C1: was a date-serial number (43751, formatted yyyy-mm-dd); now becomes the
# literal text "2019-10-13" with no special style. Force text interpretation
# via a temporary Text number format so Excel doesn't re-parse the string as
# a date, then drop back to the Normal style (which also clears the format).
$c1 = $ws.Range("C1")
$c1.NumberFormat = "@"
$c1.Value = "2019-10-13"
$c1.Style = "Normal"

# A2:A36: every attendance value flips from -7 to -13 (style untouched).
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 1).Value = -13
}

# Mark two rows present in column C.
$ws.Range("C2").Value = "Present"
$ws.Range("C35").Value = "Present"

# Move the live selection to C1 (was G5).
$ws.Range("C1").Select() | Out-Null
